$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing supplier's address (row 2, column B)
$ws.Range("B2").Value = "No 2/211 Arangala,Naula"

# New supplier: Arpico PLC (row 3)
$ws.Range("A3").Value = "Arpico PLC"
$ws.Range("B3").Value = "No 3/14,Navinna,Colombo 3"
$ws.Range("C3").Value = 762561253
$ws.Range("D3").Value = 111459823
$ws.Range("E3").Value = "arpico@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:arpico@gmail.com")
$ws.Range("E3").Style = "Hyperlink"

# New supplier: Zloan PLC (row 4)
$ws.Range("A4").Value = "Zloan PLC"
$ws.Range("B4").Value = "No 11, Kadawatha,balummahara."
$ws.Range("C4").Value = 456259635
$ws.Range("D4").Value = 385215632
$ws.Range("E4").Value = "zloan@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:zloan@gmail.com")
$ws.Range("E4").Style = "Hyperlink"

# New supplier: mitsui cement cop (row 5)
$ws.Range("A5").Value = "mitsui cement cop"
$ws.Range("B5").Value = "No 32,Trnkomalee ,Kanthale"
$ws.Range("C5").Value = 159632575
$ws.Range("D5").Value = 253974102
$ws.Range("E5").Value = "mit@yahoomail.com"
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:mit@yahoomail.com")
$ws.Range("E5").Style = "Hyperlink"

# Update selection to match post-edit cursor position
$ws.Range("D8").Select()
